$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9
# from 2023-10-13 (45212) to 2023-10-22 (45221).
$ws.Range("C2:C9").Value = 45221
